$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# --- Header text updates (shared/rich-text strings) ---
# A8: "Volume 30   Number  43" -> "...44"
$ws.Range("A8").Value = "Volume 30   Number  44"

# C9: "Report Covering the Week  10/23/2023  Through  10/29/2023"
#     -> "...10/30/2023  Through  11/5/2023"
$ws.Range("C9").Value = "Report Covering the Week  10/30/2023  Through  11/5/2023"

# --- Weekly crime-data table updates (rows 14-30, cols C:N) ---
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = -33.333333333333
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = -16.666666666666
$ws.Range("I14").Value = 52
$ws.Range("J14").Value = 58
$ws.Range("K14").Value = -10.344827586206
$ws.Range("L14").Value = 10.63829787234
$ws.Range("M14").Value = -25.714285714285
$ws.Range("N14").Value = -75.925925925925
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = -57.142857142857
$ws.Range("F15").Value = 13
$ws.Range("G15").Value = 20
$ws.Range("H15").Value = -35
$ws.Range("I15").Value = 184
$ws.Range("J15").Value = 186
$ws.Range("K15").Value = -1.075268817204
$ws.Range("L15").Value = -1.075268817204
$ws.Range("M15").Value = 20.261437908496
$ws.Range("N15").Value = -62.295081967213
$ws.Range("C16").Value = 29
$ws.Range("D16").Value = 40
$ws.Range("E16").Value = -27.5
$ws.Range("F16").Value = 138
$ws.Range("G16").Value = 150
$ws.Range("H16").Value = -8
$ws.Range("I16").Value = 1492
$ws.Range("J16").Value = 1663
$ws.Range("K16").Value = -10.282621767889
$ws.Range("L16").Value = 25.062866722548
$ws.Range("M16").Value = -39.765845781186
$ws.Range("N16").Value = -87.437905194914
$ws.Range("C17").Value = 59
$ws.Range("D17").Value = 73
$ws.Range("E17").Value = -19.17808219178
$ws.Range("F17").Value = 250
$ws.Range("G17").Value = 270
$ws.Range("H17").Value = -7.407407407407
$ws.Range("I17").Value = 2993
$ws.Range("J17").Value = 2949
$ws.Range("K17").Value = 1.492031197015
$ws.Range("L17").Value = 12.986032465081
$ws.Range("M17").Value = 43.000477783086
$ws.Range("N17").Value = -48.155205265893
$ws.Range("D18").Value = 48
$ws.Range("E18").Value = -27.083333333333
$ws.Range("F18").Value = 121
$ws.Range("G18").Value = 193
$ws.Range("H18").Value = -37.305699481865
$ws.Range("I18").Value = 1461
$ws.Range("J18").Value = 1834
$ws.Range("K18").Value = -20.338058887677
$ws.Range("L18").Value = -1.946308724832
$ws.Range("M18").Value = -50.758341759352
$ws.Range("N18").Value = -90.727930443612
$ws.Range("C19").Value = 123
$ws.Range("D19").Value = 136
$ws.Range("E19").Value = -9.558823529411
$ws.Range("F19").Value = 484
$ws.Range("G19").Value = 522
$ws.Range("H19").Value = -7.27969348659
$ws.Range("I19").Value = 5528
$ws.Range("J19").Value = 6127
$ws.Range("K19").Value = -9.776399543006
$ws.Range("L19").Value = 28.408826945412
$ws.Range("M19").Value = 18.499464094319
$ws.Range("N19").Value = -26.76205617382
$ws.Range("C20").Value = 41
$ws.Range("D20").Value = 33
$ws.Range("E20").Value = 24.242424242424
$ws.Range("F20").Value = 154
$ws.Range("G20").Value = 144
$ws.Range("H20").Value = 6.944444444444
$ws.Range("I20").Value = 1591
$ws.Range("J20").Value = 1557
$ws.Range("K20").Value = 2.18368657675
$ws.Range("L20").Value = 43.462578899909
$ws.Range("M20").Value = -3.634161114476
$ws.Range("N20").Value = -91.95652173913
$ws.Range("C21").Value = 292
$ws.Range("D21").Value = 340
$ws.Range("E21").Value = -14.117647058823
$ws.Range("F21").Value = 1165
$ws.Range("G21").Value = 1305
$ws.Range("H21").Value = -10.727969348659
$ws.Range("I21").Value = 13301
$ws.Range("J21").Value = 14374
$ws.Range("K21").Value = -7.464867121191
$ws.Range("L21").Value = 21.149467164586
$ws.Range("M21").Value = -5.505825518613
$ws.Range("N21").Value = -78.350884617262
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = -16.666666666666
$ws.Range("F22").Value = 19
$ws.Range("G22").Value = 18
$ws.Range("H22").Value = 5.555555555555
$ws.Range("I22").Value = 161
$ws.Range("J22").Value = 167
$ws.Range("K22").Value = -3.592814371257
$ws.Range("L22").Value = 29.838709677419
$ws.Range("M22").Value = -33.744855967078
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 13
$ws.Range("E23").Value = -61.538461538461
$ws.Range("F23").Value = 30
$ws.Range("G23").Value = 48
$ws.Range("H23").Value = -37.5
$ws.Range("I23").Value = 446
$ws.Range("J23").Value = 477
$ws.Range("K23").Value = -6.49895178197
$ws.Range("L23").Value = 9.313725490196
$ws.Range("M23").Value = 53.264604810996
$ws.Range("D24").Value = 344
$ws.Range("E24").Value = -16.56976744186
$ws.Range("F24").Value = 1145
$ws.Range("G24").Value = 1311
$ws.Range("H24").Value = -12.662090007627
$ws.Range("I24").Value = 13489
$ws.Range("J24").Value = 13858
$ws.Range("K24").Value = -2.662721893491
$ws.Range("L24").Value = 34.032193958664
$ws.Range("M24").Value = 27.531436135009
$ws.Range("C25").Value = 123
$ws.Range("E25").Value = 25.510204081632
$ws.Range("F25").Value = 458
$ws.Range("G25").Value = 411
$ws.Range("H25").Value = 11.435523114355
$ws.Range("I25").Value = 5002
$ws.Range("J25").Value = 4722
$ws.Range("K25").Value = 5.929690808979
$ws.Range("L25").Value = 20.29822029822
$ws.Range("M25").Value = -13.564886815275
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = -55.555555555555
$ws.Range("F26").Value = 21
$ws.Range("H26").Value = -40
$ws.Range("I26").Value = 266
$ws.Range("J26").Value = 299
$ws.Range("K26").Value = -11.036789297658
$ws.Range("L26").Value = -8.904109589041
$ws.Range("C27").Value = 11
$ws.Range("D27").Value = 11
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 65
$ws.Range("H27").Value = -26.153846153846
$ws.Range("I27").Value = 546
$ws.Range("J27").Value = 594
$ws.Range("K27").Value = -8.080808080808
$ws.Range("L27").Value = 1.298701298701
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = -83.333333333333
$ws.Range("G28").Value = 19
$ws.Range("H28").Value = -52.631578947368
$ws.Range("I28").Value = 125
$ws.Range("J28").Value = 194
$ws.Range("K28").Value = -35.567010309278
$ws.Range("L28").Value = -28.977272727272
$ws.Range("M28").Value = -46.808510638297
$ws.Range("N28").Value = -81.831395348837
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = -80
$ws.Range("G29").Value = 18
$ws.Range("H29").Value = -66.666666666666
$ws.Range("I29").Value = 108
$ws.Range("J29").Value = 152
$ws.Range("K29").Value = -28.947368421052
$ws.Range("L29").Value = -29.870129870129
$ws.Range("M29").Value = -44.615384615384
$ws.Range("N29").Value = -82.029950083194
$ws.Range("F30").Value = 12
$ws.Range("G30").Value = 11
$ws.Range("H30").Value = 9.090909090909
$ws.Range("I30").Value = 77
$ws.Range("J30").Value = 104
$ws.Range("K30").Value = -25.961538461538
$ws.Range("L30").Value = 20.3125
